# Update database values and company name for the income statement sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company name correction (shared string) ---
$ws.Range("B5").Value = "کیمیا-ص. معدنی کیمیای زنجان گستران"

# --- Period headers (row 8): drop oldest period, shift, add newest ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Publish dates (row 9): drop oldest date, shift, add newest ---
$ws.Range("D9").Value = "1399-04-19 (13)"
$ws.Range("E9").Value = "1400-04-16 (12)"
$ws.Range("F9").Value = "1401-04-08 (9)"
$ws.Range("G9").Value = "1402-02-28 (8)"

# H9 looks like a bare date ("1402-02-28"); a plain .Value assignment would
# get auto-converted by Excel into a date serial number. Enter it as a text
# formula instead (a quoted literal always yields a string result), then
# copy/paste-values over itself to bake it down into a plain static value -
# this keeps the original cell style/number format untouched and avoids
# creating a date value or a quote-prefixed style.
$ws.Range("H9").Formula = '="1402-02-28"'
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false

# --- Financial figures: shift all yearly columns left by one year and
#     append the new (1401/12) figures in column H ---

# فروش (Sales) - row 11
$ws.Range("D11").Value = 11720
$ws.Range("E11").Value = 11732
$ws.Range("F11").Value = 18783
$ws.Range("G11").Value = 26829
$ws.Range("H11").Value = 25069

# بهای تمام شده کالای فروش رفته (COGS) - row 12
$ws.Range("D12").Value = -7554
$ws.Range("E12").Value = -7325
$ws.Range("F12").Value = -12203
$ws.Range("G12").Value = -18917
$ws.Range("H12").Value = -16395

# سود (زیان) ناخالص (Gross profit) - row 13
$ws.Range("D13").Value = 4166
$ws.Range("E13").Value = 4408
$ws.Range("F13").Value = 6579
$ws.Range("G13").Value = 7912
$ws.Range("H13").Value = 8675

# هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) - row 14
$ws.Range("D14").Value = -292
$ws.Range("E14").Value = -265
$ws.Range("F14").Value = -327
$ws.Range("G14").Value = -1201
$ws.Range("H14").Value = -1027

# هزینه کاهش ارزش دریافتنی‌ها (row 16) - D stays "-", rest are numeric
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 46
$ws.Range("G16").Value = 80
$ws.Range("H16").Value = 2616

# سود (زیان) عملیاتی (Operating profit) - row 17
$ws.Range("D17").Value = 3874
$ws.Range("E17").Value = 4156
$ws.Range("F17").Value = 6299
$ws.Range("G17").Value = 6790
$ws.Range("H17").Value = 10263

# هزینه های مالی (Financial expenses) - row 18
$ws.Range("D18").Value = -27
$ws.Range("E18").Value = -129
$ws.Range("F18").Value = -66
$ws.Range("G18").Value = -14
$ws.Range("H18").Value = -36

# خالص سایر درامدها و هزینه های غیرعملیاتی - row 19
$ws.Range("D19").Value = 231
$ws.Range("E19").Value = 682
$ws.Range("F19").Value = 658
$ws.Range("G19").Value = 1070
$ws.Range("H19").Value = 528

# سود (زیان) خالص عملیات در حال تداوم قبل از مالیات - row 20
$ws.Range("D20").Value = 4077
$ws.Range("E20").Value = 4709
$ws.Range("F20").Value = 6891
$ws.Range("G20").Value = 7846
$ws.Range("H20").Value = 10755

# سود (زیان) خالص عملیات در حال تداوم - row 22
$ws.Range("D22").Value = 4077
$ws.Range("E22").Value = 4709
$ws.Range("F22").Value = 6891
$ws.Range("G22").Value = 7846
$ws.Range("H22").Value = 10755

# سود (زیان) خالص - row 24
$ws.Range("D24").Value = 4077
$ws.Range("E24").Value = 4709
$ws.Range("F24").Value = 6891
$ws.Range("G24").Value = 7846
$ws.Range("H24").Value = 10755

# سرمایه (Capital) - row 26
$ws.Range("D26").Value = 3954
$ws.Range("E26").Value = 8185
$ws.Range("F26").Value = 4644
$ws.Range("G26").Value = 3979
$ws.Range("H26").Value = 8501
